$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains a new "id_hotel" column between the existing
# "Id_habitacion" (A) and "id_servicios" (B) columns, shifting
# id_servicios -> C, Ocupacion -> D, Tarifa -> E.
$ws.Range("B1").EntireColumn.Insert()
$ws.Cells.Item(1, 2).Value = "id_hotel"

$hotelId = 10003

# Room numbers (suffix of the old "Id_habitacion" values) and the old
# "id_servicios" values, row by row (rows 2..11).
$rooms     = 101,102,103,104,105,106,107,108,109,110
$servicios = 10101,10102,10103,10104,10105,10106,10107,10108,10109,10110

for ($i = 0; $i -lt $rooms.Length; $i++) {
    $row = $i + 2

    # A: Id_habitacion becomes "<hotelId><room>"
    $ws.Cells.Item($row, 1).Value = [double]("$hotelId$($rooms[$i])")

    # B: new id_hotel column, constant hotel id
    $ws.Cells.Item($row, 2).Value = $hotelId

    # C: id_servicios (shifted from old B) becomes "<hotelId><servicio>"
    $ws.Cells.Item($row, 3).Value = [double]("$hotelId$($servicios[$i])")
}

# Column C ("id_servicios") gets an explicit width in the authored file.
$ws.Columns.Item(3).ColumnWidth = 10.33

# Selection moves onto the newly populated id_servicios column.
$ws.Range("C2:C11").Select()
